$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'54.820.67"
$ws.Range("E2").Value = "'  -2.50%  "
$ws.Range("D3").Value = "'2.328.07"
$ws.Range("E3").Value = "'  -5.84%  "
$ws.Range("E4").Value = "'  +0.19%  "
$ws.Range("D5").Value = "'470.81"
$ws.Range("E5").Value = "'  -3.34%  "
$ws.Range("D6").Value = "'143.14"
$ws.Range("E6").Value = "'  -2.17%  "
$ws.Range("E7").Value = "'  +0.30%  "
$ws.Range("D8").Value = "'0.589"
$ws.Range("E8").Value = "'  +15.68%  "
$ws.Range("D9").Value = "'2.343.68"
$ws.Range("E9").Value = "'  -5.46%  "
$ws.Range("D10").Value = "'0.0955"
$ws.Range("E10").Value = "'  -1.02%  "
$ws.Range("E11").Value = "'  -6.64%  "
$ws.Range("D12").Value = "'0.319"
$ws.Range("E12").Value = "'  -3.52%  "
$ws.Range("E13").Value = "'  +0.83%  "
$ws.Range("D14").Value = "'2.744.70"
$ws.Range("E14").Value = "'  -5.38%  "
$ws.Range("D15").Value = "'54.906.07"
$ws.Range("E15").Value = "'  -2.29%  "
$ws.Range("D16").Value = "'19.82"
$ws.Range("E16").Value = "'  -5.64%  "
$ws.Range("D17").Value = "'0.0000128"
$ws.Range("E17").Value = "'  -4.85%  "
$ws.Range("D18").Value = "'2.323.45"
$ws.Range("E18").Value = "'  -6.48%  "
$ws.Range("D19").Value = "'4.54"
$ws.Range("E19").Value = "'  +1.00%  "
$ws.Range("D20").Value = "'313.43"
$ws.Range("E20").Value = "'  -1.11%  "
$ws.Range("D21").Value = "'9.51"
$ws.Range("E21").Value = "'  -4.99%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "'  +0.16%  "
$ws.Range("D23").Value = "'5.57"
$ws.Range("E23").Value = "'  -3.60%  "
$ws.Range("D24").Value = "'56.51"
$ws.Range("E24").Value = "'  -3.13%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "'  +0.16%  "
$ws.Range("D26").Value = "'0.390"
$ws.Range("E26").Value = "'  -4.92%  "
$ws.Range("D27").Value = "'0.152"
$ws.Range("E27").Value = "'  -6.10%  "
$ws.Range("D28").Value = "'2.448.50"
$ws.Range("E28").Value = "'  -5.34%  "
$ws.Range("D29").Value = "'7.05"
$ws.Range("E29").Value = "'  -7.14%  "
$ws.Range("E30").Value = "'  +0.18%  "
$ws.Range("D31").Value = "'0.0₃0739"
$ws.Range("E31").Value = "'  -6.05%  "
$ws.Range("D32").Value = "'145.96"
$ws.Range("E32").Value = "'  -2.64%  "
$ws.Range("D33").Value = "'18.15"
$ws.Range("E33").Value = "'  +0.11%  "
$ws.Range("D34").Value = "'1.46"
$ws.Range("E34").Value = "'  -2.81%  "
$ws.Range("D35").Value = "'5.04"
$ws.Range("E35").Value = "'  -2.62%  "
$ws.Range("D36").Value = "'3.55"
$ws.Range("E36").Value = "'  -4.52%  "
$ws.Range("D37").Value = "'1.07"
$ws.Range("E37").Value = "'  -5.53%  "
$ws.Range("D38").Value = "'0.801"
$ws.Range("E38").Value = "'  -6.67%  "
$ws.Range("B39").Value = "'FirstDigitalUSD"
$ws.Range("C39").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "'  +0.71%  "
$ws.Range("B40").Value = "'OKB"
$ws.Range("C40").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'33.41"
$ws.Range("E40").Value = "'  -1.99%  "
$ws.Range("D41").Value = "'0.0975"
$ws.Range("E41").Value = "'  +5.44%  "
$ws.Range("D42").Value = "'1.31"
$ws.Range("E42").Value = "'  -0.59%  "
$ws.Range("D43").Value = "'3.37"
$ws.Range("E43").Value = "'  -4.13%  "
$ws.Range("D44").Value = "'0.575"
$ws.Range("E44").Value = "'  -5.20%  "
$ws.Range("D45").Value = "'0.0515"
$ws.Range("E45").Value = "'  -7.19%  "
$ws.Range("E46").Value = "'  -0.35%  "
$ws.Range("D47").Value = "'248.05"
$ws.Range("E47").Value = "'  -3.86%  "
$ws.Range("D48").Value = "'0.0219"
$ws.Range("E48").Value = "'  -3.84%  "
$ws.Range("D49").Value = "'4.29"
$ws.Range("E49").Value = "'  -9.77%  "
$ws.Range("D50").Value = "'16.49"
$ws.Range("E50").Value = "'  -5.47%  "
$ws.Range("D51").Value = "'1.765.27"
$ws.Range("E51").Value = "'  -5.36%  "
